# Apply "Finished added class labels to excel, cleaned data":
#  - fix two mis-assigned Class codes (K70, K83)
#  - backfill the "Class" (column K) labels for rows 90-227
#  - update the sheet selection / zoom to match the author's final view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to two already-labelled rows -----------------------------
$ws.Cells.Item(70, 11).Value = "G"
$ws.Cells.Item(83, 11).Value = "E"

# --- Backfill column K ("Class") for rows 90 through 227 -------------------
$classByRow = @{
    90 = "G"; 91 = "F"; 92 = "F"; 93 = "O"; 94 = "F"; 95 = "E"; 96 = "F"; 97 = "O"
    98 = "F"; 99 = "F"; 100 = "F"; 101 = "F"; 102 = "O"; 103 = "F"; 104 = "F"; 105 = "O"
    106 = "O"; 107 = "O"; 108 = "E"; 109 = "F"; 110 = "F"; 111 = "F"; 112 = "O"; 113 = "O"
    114 = "O"; 115 = "F"; 116 = "O"; 117 = "O"; 118 = "G"; 119 = "E"; 120 = "X"; 121 = "F"
    122 = "F"; 123 = "O"; 124 = "X"; 125 = "O"; 126 = "O"; 127 = "O"; 128 = "X"; 129 = "O"
    130 = "O"; 131 = "F"; 132 = "O"; 133 = "O"; 134 = "O"; 135 = "O"; 136 = "F"; 137 = "F"
    138 = "E"; 139 = "E"; 140 = "E"; 141 = "G"; 142 = "F"; 143 = "X"; 144 = "O"; 145 = "O"
    146 = "O"; 147 = "F"; 148 = "F"; 149 = "O"; 150 = "F"; 151 = "G"; 152 = "F"; 153 = "O"
    154 = "O"; 155 = "F"; 156 = "F"; 157 = "F"; 158 = "E"; 159 = "O"; 160 = "F"; 161 = "O"
    162 = "F"; 163 = "F"; 164 = "G"; 165 = "F"; 166 = "F"; 167 = "E"; 168 = "E"; 169 = "O"
    170 = "F"; 171 = "F"; 172 = "F"; 173 = "F"; 174 = "F"; 175 = "O"; 176 = "O"; 177 = "F"
    178 = "O"; 179 = "F"; 180 = "F"; 181 = "O"; 182 = "O"; 183 = "X"; 184 = "O"; 185 = "F"
    186 = "O"; 187 = "F"; 188 = "O"; 189 = "G"; 190 = "F"; 191 = "F"; 192 = "F"; 193 = "F"
    194 = "F"; 195 = "E"; 196 = "E"; 197 = "F"; 198 = "F"; 199 = "F"; 200 = "F"; 201 = "F"
    202 = "G"; 203 = "F"; 204 = "F"; 205 = "F"; 206 = "F"; 207 = "F"; 208 = "F"; 209 = "F"
    210 = "O"; 211 = "F"; 212 = "E"; 213 = "F"; 214 = "F"; 215 = "F"; 216 = "E"; 217 = "F"
    218 = "E"; 219 = "G"; 220 = "O"; 221 = "O"; 222 = "F"; 223 = "F"; 224 = "F"; 225 = "F"
    226 = "O"; 227 = "F"
}

foreach ($row in $classByRow.Keys) {
    $ws.Cells.Item($row, 11).Value = $classByRow[$row]
}

# --- Final view state: zoomed to 150%, scrolled near the bottom, F6 active -
$win = $excel.ActiveWindow
$win.Zoom = 150
$win.ScrollRow = 183
$win.ScrollColumn = 1
$ws.Range("F6").Select()
